# Update countries & provincias Spain
# Applies the data refresh captured by the scraper re-run (28 May 2020, 02:05)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text (row 1)
$ws.Range("A1").Value = "Datos actualizados a 28 de Mayo de 2020 a las 02:05"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1745689
$ws.Range("C4").Value = 20432
$ws.Range("D4").Value = 488465
$ws.Range("E4").Value = 1155129
$ws.Range("G4").Value = 1523
$ws.Range("H4").Value = 102095

# Row 47 - Argentina
$ws.Range("B47").Value = 13933
$ws.Range("C47").Value = 705
$ws.Range("E47").Value = 9084
$ws.Range("G47").Value = 16
$ws.Range("H47").Value = 500

# Row 165 - Guyana
$ws.Range("D165").Value = 67
$ws.Range("E165").Value = 61

# Rows 208/209 - Surinam overtakes Islas Turcas y Caicos in total cases,
# so the two countries swap places in the sorted list.
$ws.Range("A208").Value = "Surinam"
$ws.Range("C208").Value = 1
$ws.Range("D208").Value = 9
$ws.Range("E208").Value = 2

$ws.Range("A209").Value = "Islas Turcas y Caicos"
$ws.Range("B209").Value = 12
$ws.Range("D209").Value = 10
